$d = $word.ActiveDocument
$W = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# --- Edit 1 -------------------------------------------------------------
# Split the "Week 1: ..." paragraph into two paragraphs right after "Week 1: "
# and mark "trello" with spell-check proofErr tags in the new paragraph.
$week1Old = "Week 1: mostly revising things from textbook, leaning to use new software like trello and fill project documentation. Learned a lot about prototypes and prototype values "
$r1 = $d.Content
$r1.Find.Execute($week1Old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1.Text = "Week 1: " + [char]13 + "mostly revising things from textbook, leaning to use new software like trello and fill project documentation. Learned a lot about prototypes and prototype values "

# Re-locate the freshly created second paragraph by its full text and replace
# its content via InsertXML so we can interleave <w:proofErr/> markers.
$week1NewFull = "mostly revising things from textbook, leaning to use new software like trello and fill project documentation. Learned a lot about prototypes and prototype values "
$fullText1 = $d.Content.Text
$start1 = $fullText1.IndexOf($week1NewFull)
$end1 = $start1 + $week1NewFull.Length
$para1Range = $d.Range($start1, $end1)
$xml1 = '<w:p xmlns:w="' + $W + '">' + `
  '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">mostly revising things from textbook, leaning to use new software like </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>trello</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> and fill project documentation. Learned a lot about prototypes and prototype values </w:t></w:r>' + `
  '</w:p>'
$para1Range.InsertXML($xml1)

# --- Edit 2 -------------------------------------------------------------
# Within the Week 3 paragraph, split the sentence around "tho" and mark it
# with spell-check proofErr tags, keeping it as a single paragraph.
$week3Old = "Learned a lot about how programming languages work while working on chapter 12 even tho I had to look up some things. Chapter 14 was uneasy, a lot of info about HTML structure I did not know before"
$fullText2 = $d.Content.Text
$start2 = $fullText2.IndexOf($week3Old)
$end2 = $start2 + $week3Old.Length
$para2Range = $d.Range($start2, $end2)
$xml2 = '<w:p xmlns:w="' + $W + '">' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Learned a lot about how programming languages work while working on chapter 12 even </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>tho</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> I had to look up some things. Chapter 14 was uneasy, a lot of info about HTML structure I did not know before</w:t></w:r>' + `
  '</w:p>'
$para2Range.InsertXML($xml2)

Write-Output "done"
